$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("D4").Value = 91.58
$wsGrupo.Range("M4").Value = 63.8
$wsGrupo.Range("D7").Value = "1 de 5"
$wsGrupo.Range("M7").Value = "2 de 5"

# --- Sheet "VENTA MENSUAL" ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F4").Value = 155.38
$wsMensual.Range("F7").Value = 213.24
# ColumnWidth is specified in characters; the stored OOXML column width is
# ColumnWidth + 5/6, so subtract that offset to land exactly on stored width 12.
$wsMensual.Columns.Item(6).ColumnWidth = 11.166666666666666
